$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.293.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "'3.487.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'609.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'185.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.215"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "'9.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "'4.041.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "'599.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.06%  "
$ws.Range("D16").Value = "'69.358.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "'18.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "'12.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "'3.489.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").Value = "'17.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.63%  "
$ws.Range("D23").Value = "'105.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'4.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "'5.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "'3.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").Value = "'10.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "'9.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.60%  "
$ws.Range("D29").Value = "'33.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("D30").Value = "'6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "'12.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "'3.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.77%  "
$ws.Range("D34").Value = "'63.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'3.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'518.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("D38").Value = "'0.396"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("D39").Value = "'3.596.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "'0.0₃0776"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'0.0459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").Value = "'2.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("D47").Value = "'3.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.56%  "
$ws.Range("E48").Value = "  -5.04%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'0.000244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.23%  "
$ws.Range("E51").Value = "  -8.89%  "
